$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68. This shifts the existing rows 68..196
# down to 69..197, preserving all of their data (dates, prices, etc.)
# exactly as they were - which matches the "row N after = row N-1 before"
# pattern seen throughout the diff.
$ws.Rows.Item(68).EntireRow.Insert()

# Populate the newly inserted row 68 with its own (new) data record.
$ws.Cells.Item(68, 1).Value = 8
$ws.Cells.Item(68, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(68, 3).Value = "Coquimbo"
$ws.Cells.Item(68, 4).Value = 44469
$ws.Cells.Item(68, 5).Value = 4
$ws.Cells.Item(68, 6).Value = 100114013
$ws.Cells.Item(68, 7).Value = "Zanahoria"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 600
$ws.Cells.Item(68, 11).Value = 6000
$ws.Cells.Item(68, 12).Value = 7000
$ws.Cells.Item(68, 13).Value = 6500
$ws.Cells.Item(68, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(68, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(68, 16).Value = 325
$ws.Cells.Item(68, 17).Value = 20
$ws.Cells.Item(68, 18).Value = "Hortaliza"
